# Reorder the Salida column (B2:B81) on the "Solucion" sheet to match the
# new simulated_annealing/beam_search run order, and update the per-zone
# "Tiempo" values on the "Metricas" sheet (Z1 and Z4 swap values).
$wb = $excel.ActiveWorkbook

$wsSol = $wb.Worksheets.Item("Solucion")
$salidas = @(
    "S036",
    "S057",
    "S008",
    "S064",
    "S007",
    "S058",
    "S061",
    "S016",
    "S034",
    "S048",
    "S045",
    "S021",
    "S038",
    "S055",
    "S039",
    "S069",
    "S071",
    "S011",
    "S070",
    "S075",
    "S025",
    "S018",
    "S006",
    "S062",
    "S063",
    "S074",
    "S020",
    "S054",
    "S080",
    "S001",
    "S066",
    "S010",
    "S046",
    "S037",
    "S044",
    "S049",
    "S076",
    "S035",
    "S004",
    "S030",
    "S067",
    "S041",
    "S077",
    "S040",
    "S009",
    "S017",
    "S023",
    "S012",
    "S043",
    "S052",
    "S002",
    "S042",
    "S026",
    "S078",
    "S019",
    "S072",
    "S005",
    "S073",
    "S032",
    "S059",
    "S068",
    "S060",
    "S050",
    "S053",
    "S031",
    "S056",
    "S065",
    "S033",
    "S047",
    "S003",
    "S022",
    "S014",
    "S013",
    "S029",
    "S024",
    "S015",
    "S051",
    "S028",
    "S027",
    "S079"
)

for ($i = 0; $i -lt $salidas.Length; $i++) {
    $row = $i + 2
    $wsSol.Cells.Item($row, 2).Value = $salidas[$i]
}

$wsMet = $wb.Worksheets.Item("Metricas")
$wsMet.Range("B2").Value = "373,5"
$wsMet.Range("B5").Value = "373,6"
